# Fruta / hortaliza, semanal
# Insert two new weekly-report rows (Provincia de Chacabuco) right before the
# existing "Region Metropolitana" block that starts at row 458, pushing all
# the rows below it down by two (458-477 -> 460-479).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at row 458 (shifts old rows 458..477 down to 460..479).
$ws.Rows.Item(458).Resize(2).Insert()

# New row 458: Especial / Provincia de Chacabuco
$ws.Cells.Item(458, 1).Value = 9
$ws.Cells.Item(458, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(458, 3).Value = "Metropolitana"
$ws.Cells.Item(458, 4).Value = 45106
$ws.Cells.Item(458, 5).Value = 13
$ws.Cells.Item(458, 6).Value = "Fruta"
$ws.Cells.Item(458, 7).Value = 100107
$ws.Cells.Item(458, 8).Value = "Otros"
$ws.Cells.Item(458, 9).Value = 100107011
$ws.Cells.Item(458, 10).Value = "Tuna"
$ws.Cells.Item(458, 11).Value = "Sin especificar"
$ws.Cells.Item(458, 12).Value = "Especial"
$ws.Cells.Item(458, 13).Value = 250
$ws.Cells.Item(458, 14).Value = 30000
$ws.Cells.Item(458, 15).Value = 30000
$ws.Cells.Item(458, 16).Value = 30000
$ws.Cells.Item(458, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(458, 18).Value = "Provincia de Chacabuco"
$ws.Cells.Item(458, 19).Value = 1667
$ws.Cells.Item(458, 20).Value = 18

# New row 459: Primera / Provincia de Chacabuco
$ws.Cells.Item(459, 1).Value = 9
$ws.Cells.Item(459, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(459, 3).Value = "Metropolitana"
$ws.Cells.Item(459, 4).Value = 45106
$ws.Cells.Item(459, 5).Value = 13
$ws.Cells.Item(459, 6).Value = "Fruta"
$ws.Cells.Item(459, 7).Value = 100107
$ws.Cells.Item(459, 8).Value = "Otros"
$ws.Cells.Item(459, 9).Value = 100107011
$ws.Cells.Item(459, 10).Value = "Tuna"
$ws.Cells.Item(459, 11).Value = "Sin especificar"
$ws.Cells.Item(459, 12).Value = "Primera"
$ws.Cells.Item(459, 13).Value = 180
$ws.Cells.Item(459, 14).Value = 26000
$ws.Cells.Item(459, 15).Value = 26000
$ws.Cells.Item(459, 16).Value = 26000
$ws.Cells.Item(459, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(459, 18).Value = "Provincia de Chacabuco"
$ws.Cells.Item(459, 19).Value = 1444
$ws.Cells.Item(459, 20).Value = 18
